$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update TOTAL GRANT REQUESTED value
$ws.Range("B7").Value = 8970

# Add new line item: miscellaneous expenses
$ws.Range("A23").Value = "miscellaneous expenses (meals, local transportation,..)"
$ws.Range("B23").Value = 1600
